# Pdcd1lg2-Rgmb LR-pairs sheet: refresh with new TPM-derived NATMI output.
# Existing rows 2-5 get updated values; 4 additional sending/target-cluster
# combinations (rows 6-9) are appended for the "Resolving-Mac" sender.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A-T => 1-20. Col order:
# A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
# E..T the 16 numeric NATMI metrics.
$newData = @(
    @("ECs", "Pdcd1lg2", "Rgmb", "ECs", 1, 0.3333333333333333, 0.04600433333333333, 0.138013, 0.06627080791390985, 0.06627080791390984, 3, 1, 11.188041, 33.564123, 0.1395103797998223, 0.1395103797998223, 0.5146983675110001, 4.632285307599, 0.009245465581710632, 0.00924546558171063),
    @("ECs", "Pdcd1lg2", "Rgmb", "FAPs", 1, 0.3333333333333333, 0.04600433333333333, 0.138013, 0.06627080791390985, 0.06627080791390984, 3, 1, 33.42883866666667, 100.286516, 0.4168441980730721, 0.4168441980730722, 1.537871436967556, 13.840842932708, 0.02762460178052835, 0.02762460178052835),
    @("ECs", "Pdcd1lg2", "Rgmb", "MuSCs", 1, 0.3333333333333333, 0.04600433333333333, 0.138013, 0.06627080791390985, 0.06627080791390984, 3, 1, 31.78201566666667, 95.346047, 0.3963089763847458, 0.3963089763847459, 1.462110442734555, 13.158993984611, 0.02626371604855172, 0.02626371604855173),
    @("ECs", "Pdcd1lg2", "Rgmb", "Resolving-Mac", 1, 0.3333333333333333, 0.04600433333333333, 0.138013, 0.06627080791390985, 0.06627080791390984, 3, 1, 3.796148333333333, 11.388445, 0.04733644574235969, 0.04733644574235969, 0.1746392733094444, 1.571753459785, 0.003137024503119134, 0.003137024503119134),
    @("Resolving-Mac", "Pdcd1lg2", "Rgmb", "ECs", 3, 1, 0.6481826666666667, 1.944548, 0.9337291920860902, 0.9337291920860902, 3, 1, 11.188041, 33.564123, 0.1395103797998223, 0.1395103797998223, 7.251894250156, 65.267048251404, 0.1302649142181116, 0.1302649142181117),
    @("Resolving-Mac", "Pdcd1lg2", "Rgmb", "FAPs", 3, 1, 0.6481826666666667, 1.944548, 0.9337291920860902, 0.9337291920860902, 3, 1, 33.42883866666667, 100.286516, 0.4168441980730721, 0.4168441980730722, 21.66799379052978, 195.011944114768, 0.3892195962925438, 0.3892195962925438),
    @("Resolving-Mac", "Pdcd1lg2", "Rgmb", "MuSCs", 3, 1, 0.6481826666666667, 1.944548, 0.9337291920860902, 0.9337291920860902, 3, 1, 31.78201566666667, 95.346047, 0.3963089763847458, 0.3963089763847459, 20.60055166686178, 185.404965001756, 0.3700452603361941, 0.3700452603361942),
    @("Resolving-Mac", "Pdcd1lg2", "Rgmb", "Resolving-Mac", 3, 1, 0.6481826666666667, 1.944548, 0.9337291920860902, 0.9337291920860902, 3, 1, 3.796148333333333, 11.388445, 0.04733644574235969, 0.04733644574235969, 2.460597549762222, 22.14537794786, 0.04419942123924055, 0.04419942123924056),
)

for ($i = 0; $i -lt $newData.Count; $i++) {
    $rowNum = $i + 2
    $values = $newData[$i]
    for ($col = 1; $col -le $values.Count; $col++) {
        $ws.Cells.Item($rowNum, $col).Value = $values[$col - 1]
    }
}

Write-Output "Updated $($newData.Count) data rows (A2:T$($newData.Count + 1))"
